$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.Style = "Normal"
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "257.53"
Set-TextValue "E2" "0.10%"
Set-TextValue "E3" "-1.33%"
Set-TextValue "D4" "4.689"
Set-TextValue "E4" "-9.94%"
Set-TextValue "D5" "0.05883"
Set-TextValue "E5" "-0.51%"
Set-TextValue "D6" "6.648"
Set-TextValue "E6" "-0.34%"
Set-TextValue "D7" "0.8576"
Set-TextValue "E7" "-0.93%"
Set-TextValue "D8" "0.9622"
Set-TextValue "E8" "-5.76%"
Set-TextValue "D9" "0.1410"
Set-TextValue "E9" "-0.61%"
Set-TextValue "D10" "0.03979"
Set-TextValue "E10" "10.98%"
Set-TextValue "D11" "0.07087"
Set-TextValue "E11" "-1.34%"
Set-TextValue "D12" "0.03180"
Set-TextValue "E12" "0.75%"
Set-TextValue "D13" "0.09173"
Set-TextValue "E13" "-0.62%"
Set-TextValue "D14" "0.001542"
Set-TextValue "E14" "0.36%"
Set-TextValue "D15" "0.0006063"
Set-TextValue "E15" "-94.21%"
Set-TextValue "D16" "0.006205"
Set-TextValue "E16" "4.07%"
Set-TextValue "E17" "1.09%"
Set-TextValue "D18" "3.204"
Set-TextValue "E18" "-1.83%"
Set-TextValue "E19" "-0.81%"
Set-TextValue "D20" "0.3079"
Set-TextValue "E20" "-2.30%"
Set-TextValue "E21" "-1.04%"
Set-TextValue "D22" "3.866"
Set-TextValue "E22" "9.80%"
Set-TextValue "D23" "0.04217"
Set-TextValue "E23" "1.05%"
Set-TextValue "D24" "0.001220"
Set-TextValue "E24" "0.19%"
Set-TextValue "D25" "0.004295"
Set-TextValue "E25" "-4.87%"
Set-TextValue "D27" "0.0001938"
Set-TextValue "E27" "0.02%"
Set-TextValue "D40" "0.03829"
Set-TextValue "E40" "0.07%"
Set-TextValue "D41" "0.006190"
Set-TextValue "E41" "12.37%"
Set-TextValue "E42" "0.03%"
Set-TextValue "E43" "15.82%"
Set-TextValue "E44" "7.09%"
Set-TextValue "D45" "0.00005463"
Set-TextValue "E45" "0.54%"
Set-TextValue "E46" "0.03%"
Set-TextValue "D47" "0.06003"
Set-TextValue "E47" "-44.96%"
Set-TextValue "D48" "0.1754"
Set-TextValue "E48" "7,962.92%"
Set-TextValue "E49" "0.03%"
Set-TextValue "E50" "0.03%"
